# Regenerate orders with updated distance/size codes:
#   Distance: D64 -> D69, D80 -> D86, D51 -> D55
#   Size:     S30 -> S31   (S25 and S20 are left unchanged)
#
# These codes show up inside many different text values on the sheet
# (condition names, left/right filenames, the Distance/Size lookup
# columns, etc.), always as a literal substring of a larger token such
# as "Face03_D80_S25" or "Face14_D80_S30_l.png". We walk every used
# cell, and for any text cell whose value contains one of the old
# tokens we rewrite it with the new token(s) substituted in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$changed = 0

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal -replace 'D64', 'D69'
            $newVal = $newVal -replace 'D80', 'D86'
            $newVal = $newVal -replace 'D51', 'D55'
            $newVal = $newVal -replace 'S30', 'S31'

            if ($newVal -ne $val) {
                $cell.Value = $newVal
                $changed = $changed + 1
            }
        }
    }
}

Write-Output "Updated $changed cell(s) with new distance/size codes."
